# Regenerate save_data: replace the "Strike#" values in column G ("K") with
# freshly calculated K values (s_vals), per regenerated std/mean calc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 1
    7  = 2
    8  = 0
    9  = 2
    10 = 0
    11 = 1
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    18 = 3
    19 = 0
    20 = 2
    21 = 1
    22 = 3
    23 = 0
    24 = 0
    25 = 3
    26 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
